$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data values (as text, matching original inline-string cells)
# and reset each touched cell back to the default (unstyled) format so only the
# value changes - matching the target diff which shows no style/format changes.

# Row 2
$ws.Range("D2").Value = "'295.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-3.37%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'10"
$ws.Range("G2").Style = "Normal"
# Row 3
$ws.Range("D3").Value = "'31.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-3.54%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'10"
$ws.Range("G3").Style = "Normal"
# Row 4
$ws.Range("D4").Value = "'5.109"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.77%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'10"
$ws.Range("G4").Style = "Normal"
# Row 5
$ws.Range("D5").Value = "'0.07371"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.95%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'10"
$ws.Range("G5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = "'7.693"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.50%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'10"
$ws.Range("G6").Style = "Normal"
# Row 7
$ws.Range("B7").Value = "'FTXToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'1.658"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'12.06%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'10"
$ws.Range("G7").Style = "Normal"
# Row 8
$ws.Range("B8").Value = "'GateToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'3.748"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.19%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'10"
$ws.Range("G8").Style = "Normal"
# Row 9
$ws.Range("D9").Value = "'0.9161"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.66%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'10"
$ws.Range("G9").Style = "Normal"
# Row 10
$ws.Range("D10").Value = "'0.1669"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.21%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'10"
$ws.Range("G10").Style = "Normal"
# Row 11
$ws.Range("D11").Value = "'0.07187"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-3.29%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'10"
$ws.Range("G11").Style = "Normal"
# Row 12
$ws.Range("D12").Value = "'0.07957"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.41%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'10"
$ws.Range("G12").Style = "Normal"
# Row 13
$ws.Range("D13").Value = "'0.02992"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.24%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'10"
$ws.Range("G13").Style = "Normal"
# Row 14
$ws.Range("D14").Value = "'0.09890"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.45%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'10"
$ws.Range("G14").Style = "Normal"
# Row 15
$ws.Range("D15").Value = "'0.001516"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.66%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'10"
$ws.Range("G15").Style = "Normal"
# Row 16
$ws.Range("D16").Value = "'0.006188"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-4.97%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'10"
$ws.Range("G16").Style = "Normal"
# Row 17
$ws.Range("E17").Value = "'-1.37%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'10"
$ws.Range("G17").Style = "Normal"
# Row 18
$ws.Range("D18").Value = "'2.229"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.02%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'10"
$ws.Range("G18").Style = "Normal"
# Row 19
$ws.Range("D19").Value = "'0.3272"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.84%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'10"
$ws.Range("G19").Style = "Normal"
# Row 20
$ws.Range("D20").Value = "'0.1346"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.63%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'10"
$ws.Range("G20").Style = "Normal"
# Row 21
$ws.Range("D21").Value = "'4.549"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'6.50%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'10"
$ws.Range("G21").Style = "Normal"
# Row 22
$ws.Range("D22").Value = "'0.04619"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.48%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'10"
$ws.Range("G22").Style = "Normal"
# Row 23
$ws.Range("D23").Value = "'0.1548"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-5.14%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'10"
$ws.Range("G23").Style = "Normal"
# Row 24
$ws.Range("D24").Value = "'0.001216"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.75%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'10"
$ws.Range("G24").Style = "Normal"
# Row 25
$ws.Range("D25").Value = "'0.004423"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.07%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'10"
$ws.Range("G25").Style = "Normal"
# Row 26
$ws.Range("D26").Value = "'0.0001297"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.71%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'10"
$ws.Range("G26").Style = "Normal"
# Row 27
$ws.Range("D27").Value = "'0.0001872"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'7.11%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'10"
$ws.Range("G27").Style = "Normal"
# Row 28
$ws.Range("G28").Value = "'10"
$ws.Range("G28").Style = "Normal"
# Row 29
$ws.Range("G29").Value = "'10"
$ws.Range("G29").Style = "Normal"
# Row 30
$ws.Range("G30").Value = "'10"
$ws.Range("G30").Style = "Normal"
# Row 31
$ws.Range("G31").Value = "'10"
$ws.Range("G31").Style = "Normal"
# Row 32
$ws.Range("G32").Value = "'10"
$ws.Range("G32").Style = "Normal"
# Row 33
$ws.Range("G33").Value = "'10"
$ws.Range("G33").Style = "Normal"
# Row 34
$ws.Range("G34").Value = "'10"
$ws.Range("G34").Style = "Normal"
# Row 35
$ws.Range("G35").Value = "'10"
$ws.Range("G35").Style = "Normal"
# Row 36
$ws.Range("G36").Value = "'10"
$ws.Range("G36").Style = "Normal"
# Row 37
$ws.Range("G37").Value = "'10"
$ws.Range("G37").Style = "Normal"
# Row 38
$ws.Range("G38").Value = "'10"
$ws.Range("G38").Style = "Normal"
# Row 39
$ws.Range("D39").Value = "'0.01681"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.53%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'10"
$ws.Range("G39").Style = "Normal"
# Row 40
$ws.Range("D40").Value = "'0.04414"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-2.36%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'10"
$ws.Range("G40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = "'0.007076"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.76%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'10"
$ws.Range("G41").Style = "Normal"
# Row 42
$ws.Range("D42").Value = "'0.1327"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.27%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'10"
$ws.Range("G42").Style = "Normal"
# Row 43
$ws.Range("D43").Value = "'0.002105"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-11.60%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'10"
$ws.Range("G43").Style = "Normal"
# Row 44
$ws.Range("D44").Value = "'0.01115"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-13.66%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'10"
$ws.Range("G44").Style = "Normal"
# Row 45
$ws.Range("D45").Value = "'0.00005989"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.73%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'10"
$ws.Range("G45").Style = "Normal"
# Row 46
$ws.Range("G46").Value = "'10"
$ws.Range("G46").Style = "Normal"
# Row 47
$ws.Range("D47").Value = "'0.01098"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-15.98%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'10"
$ws.Range("G47").Style = "Normal"
# Row 48
$ws.Range("G48").Value = "'10"
$ws.Range("G48").Style = "Normal"
# Row 49
$ws.Range("G49").Value = "'10"
$ws.Range("G49").Style = "Normal"
# Row 50
$ws.Range("G50").Value = "'10"
$ws.Range("G50").Style = "Normal"
# Row 51
$ws.Range("G51").Value = "'10"
$ws.Range("G51").Style = "Normal"

Write-Host "Applied all changes"
